$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells we touch so numeric-looking strings
# (e.g. "1.00", "0.160") are preserved exactly as text, matching the
# original inline-string cell type, instead of being coerced to numbers.
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '98.526.03'
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -0.26%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.466.62'
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  +4.54%  '
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '260.43'
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +1.56%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '670.98'
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  +7.34%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.55'
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  +7.85%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.457'
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +13.28%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '1.11'
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  +20.37%  '
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -0.06%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '3.463.27'
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +4.52%  '
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +9.17%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '43.27'
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  +9.65%  '
$c = $ws.Range("B14")
$c.NumberFormat = "@"
$c.Value = 'Toncoin'
$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.29'
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  +14.49%  '
$c = $ws.Range("B15")
$c.NumberFormat = "@"
$c.Value = 'ShibaInu'
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.0000271'
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +8.35%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '98.058.07'
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -0.39%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '4.114.51'
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  +4.68%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '8.81'
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +38.65%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '3.458.86'
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  +4.75%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.84'
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +15.69%  '
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +73.57%  '
$c = $ws.Range("B22")
$c.NumberFormat = "@"
$c.Value = 'Uniswap'
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '11.25'
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  +18.59%  '
$c = $ws.Range("B23")
$c.NumberFormat = "@"
$c.Value = 'SuiNetwork'
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '3.58'
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  +2.15%  '
$c = $ws.Range("B24")
$c.NumberFormat = "@"
$c.Value = 'BitcoinCash'
$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '523.44'
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  +7.73%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.0000216'
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +5.13%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '6.47'
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  +14.63%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '103.11'
$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  +16.33%  '
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +8.66%  '
$c = $ws.Range("B29")
$c.NumberFormat = "@"
$c.Value = 'Hedera'
$c = $ws.Range("C29")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.160'
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  +22.22%  '
$c = $ws.Range("B30")
$c.NumberFormat = "@"
$c.Value = 'InternetComputer(DFINITY)'
$c = $ws.Range("C30")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '11.82'
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  +15.68%  '
$c = $ws.Range("B31")
$c.NumberFormat = "@"
$c.Value = 'Cronos'
$c = $ws.Range("C31")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.197'
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +4.51%  '
$c = $ws.Range("B32")
$c.NumberFormat = "@"
$c.Value = 'Dai'
$c = $ws.Range("C32")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -0.26%  '
$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = 'PolygonEcosystemToken'
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.603'
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +30.01%  '
$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = 'Binance-PegBSC-USD'
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -0.10%  '
$c = $ws.Range("B35")
$c.NumberFormat = "@"
$c.Value = 'EthereumClassic'
$c = $ws.Range("C35")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '30.57'
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +9.60%  '
$c = $ws.Range("B36")
$c.NumberFormat = "@"
$c.Value = 'PancakeSwap'
$c = $ws.Range("C36")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.25'
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  +15.59%  '
$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = 'RenderToken'
$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '8.11'
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  +12.02%  '
$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = 'Kaspa'
$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.161'
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  +8.60%  '
$c = $ws.Range("B39")
$c.NumberFormat = "@"
$c.Value = 'Bittensor'
$c = $ws.Range("C39")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '535.18'
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +8.70%  '
$c = $ws.Range("B40")
$c.NumberFormat = "@"
$c.Value = 'Fetch.AI'
$c = $ws.Range("C40")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '1.43'
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  +15.22%  '
$c = $ws.Range("B41")
$c.NumberFormat = "@"
$c.Value = 'WhiteBITCoin'
$c = $ws.Range("C41")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '24.77'
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -0.31%  '
$c = $ws.Range("B42")
$c.NumberFormat = "@"
$c.Value = 'VeChain'
$c = $ws.Range("C42")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.0450'
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  +36.32%  '
$c = $ws.Range("B43")
$c.NumberFormat = "@"
$c.Value = 'ARBITRUM'
$c = $ws.Range("C43")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.872'
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +9.54%  '
$c = $ws.Range("B44")
$c.NumberFormat = "@"
$c.Value = 'MantraDAO'
$c = $ws.Range("C44")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '3.73'
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  +2.58%  '
$c = $ws.Range("B45")
$c.NumberFormat = "@"
$c.Value = 'dogwifhat'
$c = $ws.Range("C45")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.45'
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +9.55%  '
$c = $ws.Range("B46")
$c.NumberFormat = "@"
$c.Value = 'Cosmos'
$c = $ws.Range("C46")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '8.58'
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  +14.25%  '
$c = $ws.Range("B47")
$c.NumberFormat = "@"
$c.Value = 'Filecoin'
$c = $ws.Range("C47")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '5.40'
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +14.90%  '
$c = $ws.Range("B48")
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'
$c = $ws.Range("C48")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '1.61'
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  +18.06%  '
$c = $ws.Range("B49")
$c.NumberFormat = "@"
$c.Value = 'USDe'
$c = $ws.Range("C49")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  +0.03%  '
$c = $ws.Range("B50")
$c.NumberFormat = "@"
$c.Value = 'Stacks'
$c = $ws.Range("C50")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '2.12'
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +8.97%  '
$c = $ws.Range("B51")
$c.NumberFormat = "@"
$c.Value = 'OKB'
$c = $ws.Range("C51")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '52.17'
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  +15.23%  '
